# Apply scheduled market-data refresh to Leve profit columns (H-N) across all job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 3222.125
$ws.Range("I19").Value = 795.6667
$ws.Range("K19").Value = 795.6667
$ws.Range("M19").Value = -620.6667
$ws.Range("H40").Value = 4627.136
$ws.Range("J40").Value = 2732.6667
$ws.Range("L40").Value = 2732.6667
$ws.Range("N40").Value = -3082.6667
$ws.Range("H113").Value = 1958.8889
$ws.Range("I113").Value = 1771.3334
$ws.Range("K113").Value = 1771.3334
$ws.Range("M113").Value = 1482.6666
$ws.Range("H125").Value = 2183.75
$ws.Range("I125").Value = 1396.8572
$ws.Range("K125").Value = 12571.7148
$ws.Range("M125").Value = -10111.7148
$ws.Range("H138").Value = 1324.125
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 1324.125
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 3972.375
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -14252.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4423.3936
$ws.Range("I32").Value = 3104.7856
$ws.Range("K32").Value = 3104.7856
$ws.Range("M32").Value = -2817.7856
$ws.Range("H45").Value = 19648.75
$ws.Range("I45").Value = 33837.25
$ws.Range("K45").Value = 33837.25
$ws.Range("M45").Value = -33460.25
$ws.Range("H63").Value = 2895.58
$ws.Range("I63").Value = 1766.5
$ws.Range("J63").Value = 3110.6428
$ws.Range("K63").Value = 1766.5
$ws.Range("L63").Value = 3110.6428
$ws.Range("M63").Value = -1080.5
$ws.Range("N63").Value = -4482.6428
$ws.Range("H66").Value = 2895.58
$ws.Range("I66").Value = 1766.5
$ws.Range("J66").Value = 3110.6428
$ws.Range("K66").Value = 8832.5
$ws.Range("L66").Value = 15553.214
$ws.Range("M66").Value = -5400.5
$ws.Range("N66").Value = -22417.214
$ws.Range("H74").Value = 9422.5
$ws.Range("I74").Value = 9975.4375
$ws.Range("K74").Value = 9975.4375
$ws.Range("M74").Value = -9101.4375
$ws.Range("H77").Value = 9422.5
$ws.Range("I77").Value = 9975.4375
$ws.Range("K77").Value = 49877.1875
$ws.Range("M77").Value = -45509.1875
$ws.Range("H88").Value = 999.6667
$ws.Range("J88").Value = 999.6667
$ws.Range("L88").Value = 999.6667
$ws.Range("N88").Value = -1811.6667
$ws.Range("H91").Value = 999.6667
$ws.Range("J91").Value = 999.6667
$ws.Range("L91").Value = 999.6667
$ws.Range("N91").Value = -3807.6667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3373.1538
$ws.Range("I99").Value = 3076.3333
$ws.Range("K99").Value = 3076.3333
$ws.Range("M99").Value = -1578.3333
$ws.Range("H105").Value = 4249.645
$ws.Range("I105").Value = 2928
$ws.Range("J105").Value = 11122.2
$ws.Range("K105").Value = 2928
$ws.Range("L105").Value = 11122.2
$ws.Range("M105").Value = -1181
$ws.Range("N105").Value = -14616.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3109.2727
$ws.Range("I16").Value = 2700.6
$ws.Range("J16").Value = 3449.8333
$ws.Range("K16").Value = 2700.6
$ws.Range("L16").Value = 3449.8333
$ws.Range("M16").Value = -2413.6
$ws.Range("N16").Value = -4023.8333
$ws.Range("H31").Value = 3563.5881
$ws.Range("I31").Value = 3437.889
$ws.Range("K31").Value = 3437.889
$ws.Range("M31").Value = -3142.889
$ws.Range("H34").Value = 3563.5881
$ws.Range("I34").Value = 3437.889
$ws.Range("K34").Value = 3437.889
$ws.Range("M34").Value = -3235.889
$ws.Range("H59").Value = 199166.67
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()
$ws.Range("H113").Value = 3109.2727
$ws.Range("I113").Value = 2700.6
$ws.Range("J113").Value = 3449.8333
$ws.Range("K113").Value = 2700.6
$ws.Range("L113").Value = 3449.8333
$ws.Range("M113").Value = -530.5999999999999
$ws.Range("N113").Value = -7789.8333
$ws.Range("H122").Value = 1892.5714
$ws.Range("I122").Value = 1902.8667
$ws.Range("J122").Value = 1886.8518
$ws.Range("K122").Value = 5708.6001
$ws.Range("L122").Value = 5660.555399999999
$ws.Range("M122").Value = -3258.6001
$ws.Range("N122").Value = -10560.5554

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 1187
$ws.Range("I51").Value = 1072.875
$ws.Range("K51").Value = 3218.625
$ws.Range("M51").Value = -2758.625
$ws.Range("H118").Value = 979.5
$ws.Range("I118").Value = 979.5
$ws.Range("K118").Value = 2938.5
$ws.Range("M118").Value = -1695.5
$ws.Range("H121").Value = 404.09836
$ws.Range("J121").Value = 407.44067
$ws.Range("L121").Value = 1222.32201
$ws.Range("N121").Value = -3842.32201
$ws.Range("H131").Value = 895769.9399999999
$ws.Range("I131").Value = 3677468
$ws.Range("K131").Value = 11032404
$ws.Range("M131").Value = -11027364
$ws.Range("H141").Value = 1775.5
$ws.Range("I141").Value = 1775.5
$ws.Range("K141").Value = 5326.5
$ws.Range("M141").Value = -146.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5070.6665
$ws.Range("I70").Value = 4686.7144
$ws.Range("J70").Value = 5608.2
$ws.Range("K70").Value = 4686.7144
$ws.Range("L70").Value = 5608.2
$ws.Range("M70").Value = -4416.7144
$ws.Range("N70").Value = -6148.2
$ws.Range("H73").Value = 5070.6665
$ws.Range("I73").Value = 4686.7144
$ws.Range("J73").Value = 5608.2
$ws.Range("K73").Value = 4686.7144
$ws.Range("L73").Value = 5608.2
$ws.Range("M73").Value = -3750.7144
$ws.Range("N73").Value = -7480.2
$ws.Range("H126").Value = 7726.125
$ws.Range("I126").Value = 8452.200000000001
$ws.Range("K126").Value = 25356.6
$ws.Range("M126").Value = -22886.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 20837670
$ws.Range("I16").Value = 25003606
$ws.Range("K16").Value = 25003606
$ws.Range("M16").Value = -25003436
$ws.Range("H122").Value = 6328.122
$ws.Range("I122").Value = 6870.154
$ws.Range("J122").Value = 5388.6
$ws.Range("K122").Value = 20610.462
$ws.Range("L122").Value = 16165.8
$ws.Range("M122").Value = -18160.462
$ws.Range("N122").Value = -21065.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 411.41934
$ws.Range("I113").Value = 242.27272
$ws.Range("K113").Value = 726.81816
$ws.Range("M113").Value = 1443.18184
$ws.Range("H122").Value = 11205.05
$ws.Range("I122").Value = 7067.067
$ws.Range("J122").Value = 23619
$ws.Range("K122").Value = 21201.201
$ws.Range("L122").Value = 70857
$ws.Range("M122").Value = -18751.201
$ws.Range("N122").Value = -75757
